# donor-v0.yaml update: the "blood_type_suggested" / "Submitter Suggestion"
# column + list option were dropped from the donor schema. Remove the
# corresponding column from the TSV-export sheet and the corresponding
# list entry from the "blood_type list" sheet, then repair the header
# comments and the blood_type dropdown validation that referenced them.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Drop the "blood_type_suggested" column (D). This shifts every column
#    after it one place to the left, including the header cells and the
#    per-column data validations (their sqref ranges move automatically).
$ws1.Range("D1").EntireColumn.Delete()

# 2. Drop the "Submitter Suggestion" entry from the blood_type list sheet.
$ws2.Range("A5").EntireRow.Delete()

# 3. Comments are anchored to fixed cells and do not slide together with
#    the column delete, so re-point each header's comment at the text
#    that now belongs under it (i.e. shift the comment text left by one
#    column, same as the headers already did).
[void]$ws1.Range("D1").Comment.Text('An individual''s weight in kilograms divided by the square of the height in meters.')
[void]$ws1.Range("E1").Comment.Text('The circumstance or condition that caused death.')
[void]$ws1.Range("F1").Comment.Text('Unit for height measurement.')
[void]$ws1.Range("G1").Comment.Text('The vertical measurement or distance from the base to the top of a subject or participant.')
[void]$ws1.Range("H1").Comment.Text('The Kidney Donor Profle Index (KDPI) is a numerical measure that combines ten donor factors, including clinical parameters and demographics, to summarize into a single number the quality of deceased donor kidneys relative to other recovered kidneys. The KDPI is derived by frst calculating the Kidney Donor Risk Index (KDRI) for a deceased donor. Kidneys from a donor with a KDPI of 90%, for example, have a KDRI (which indicates relative risk of graft failure) greater than 90% of recovered kidneys. The KDPI is simply a mapping of the KDRI from a relative risk scale to a cumulative percentage scale. The reference population used for this mapping is all deceased donors in the United States with a kidney recovered for the purpose of transplantation in the prior calendar year. Lower KDPI values are associated with increased donor quality and expected longevity. https://optn.transplant.hrsa.gov/media/1512/guide_to_calculating_interpreting_kdpi.pdf
')
[void]$ws1.Range("I1").Comment.Text('Mechanism of injury may be, for example: fall, impact (eg: auto accident), weapon (eg: firearm), etc.')
[void]$ws1.Range("J1").Comment.Text('A record of a patient''s background regarding health and the occurrence of disease events of the individual.')
[void]$ws1.Range("K1").Comment.Text('A grouping of humans based on shared physical characteristics or social/ethnic identity generally viewed as distinct.')
[void]$ws1.Range("L1").Comment.Text('Biological sex at birth: male or female or other.')
[void]$ws1.Range("M1").Comment.Text('Unit for weight measurement.')
[void]$ws1.Range("N1").Comment.Text('A measurement that describes the vertical force exerted by a mass of the patient as a result of gravity.')

# The trailing comment (previously on O1, now off the end of the sheet)
# is no longer needed.
[void]$ws1.Range("O1").Comment.Delete()

# 4. The blood_type dropdown validation on column C still points at the
#    old 5-row list range / mentions the removed option in its error text;
#    fix both up to match the now 4-row list.
$dv = $ws1.Range("C2:C1048576").Validation
$dv.Modify(3, 1, 1, '''blood_type list''!$A$1:$A$4')
$dv.ErrorMessage = "Value must be one of: A / B / AB / O."
